$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like E1 (bold header style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

# Data cells F2:F66 - time_taken timestamps
$timestamps = @(
    "2021-10-05 13:41:18.157918",
    "2021-10-05 13:41:18.157932",
    "2021-10-05 13:41:18.157936",
    "2021-10-05 13:41:18.157939",
    "2021-10-05 13:41:18.157942",
    "2021-10-05 13:41:18.157945",
    "2021-10-05 13:41:18.157948",
    "2021-10-05 13:41:18.157951",
    "2021-10-05 13:41:18.157955",
    "2021-10-05 13:41:18.157958",
    "2021-10-05 13:41:18.157961",
    "2021-10-05 13:41:18.157964",
    "2021-10-05 13:41:18.157967",
    "2021-10-05 13:41:18.157970",
    "2021-10-05 13:41:18.157972",
    "2021-10-05 13:41:18.157975",
    "2021-10-05 13:41:18.157979",
    "2021-10-05 13:41:18.157982",
    "2021-10-05 13:41:18.157985",
    "2021-10-05 13:41:18.157988",
    "2021-10-05 13:41:18.157991",
    "2021-10-05 13:41:18.157994",
    "2021-10-05 13:41:18.157997",
    "2021-10-05 13:41:18.158000",
    "2021-10-05 13:41:18.158003",
    "2021-10-05 13:41:18.158006",
    "2021-10-05 13:41:18.158014",
    "2021-10-05 13:41:18.158020",
    "2021-10-05 13:41:18.158025",
    "2021-10-05 13:41:18.158030",
    "2021-10-05 13:41:18.158035",
    "2021-10-05 13:41:18.158041",
    "2021-10-05 13:41:18.158046",
    "2021-10-05 13:41:18.158052",
    "2021-10-05 13:41:18.158055",
    "2021-10-05 13:41:18.158058",
    "2021-10-05 13:41:18.158061",
    "2021-10-05 13:41:18.158064",
    "2021-10-05 13:41:18.158067",
    "2021-10-05 13:41:18.158070",
    "2021-10-05 13:41:18.158073",
    "2021-10-05 13:41:18.158077",
    "2021-10-05 13:41:18.158080",
    "2021-10-05 13:41:18.158082",
    "2021-10-05 13:41:18.158085",
    "2021-10-05 13:41:18.158088",
    "2021-10-05 13:41:18.158091",
    "2021-10-05 13:41:18.158094",
    "2021-10-05 13:41:18.158097",
    "2021-10-05 13:41:18.158100",
    "2021-10-05 13:41:18.158103",
    "2021-10-05 13:41:18.158106",
    "2021-10-05 13:41:18.158109",
    "2021-10-05 13:41:18.158112",
    "2021-10-05 13:41:18.158115",
    "2021-10-05 13:41:18.158118",
    "2021-10-05 13:41:18.158121",
    "2021-10-05 13:41:18.158124",
    "2021-10-05 13:41:18.158127",
    "2021-10-05 13:41:18.158130",
    "2021-10-05 13:41:18.158133",
    "2021-10-05 13:41:18.158136",
    "2021-10-05 13:41:18.158139",
    "2021-10-05 13:41:18.158142",
    "2021-10-05 13:41:18.158147"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
